$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3378
$ws.Range("E2").Value = 174
$ws.Range("F2").Value = 174
$ws.Range("G2").Value = 322
$ws.Range("H2").Value = 264
$ws.Range("I2").Value = 264
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 3057
$ws.Range("L2").Value = 1207
$ws.Range("M2").Value = 1850
$ws.Range("N2").Value = 1850
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 210
$ws.Range("Q2").Value = 155
$ws.Range("R2").Value = -101
$ws.Range("S2").Value = -42
$ws.Range("T2").Value = 222
$ws.Range("U2").Value = -67
$ws.Range("V2").Value = 750
$ws.Range("W2").Value = 5.14
$ws.Range("X2").Value = 7.81
$ws.Range("Y2").Value = 15.28
$ws.Range("Z2").Value = 8.640000000000001
$ws.Range("AA2").Value = 65.26000000000001
$ws.Range("AB2").Value = 780
$ws.Range("AC2").Value = 628
$ws.Range("AD2").Value = 8.779999999999999
$ws.Range("AE2").Value = 4480
$ws.Range("AG2").Value = 70
$ws.Range("AH2").Value = 1.27
$ws.Range("AI2").Value = 11.03
$ws.Range("AJ2").Value = 36720000

# Row 3
$ws.Range("D3").Value = 3765
$ws.Range("E3").Value = 222
$ws.Range("F3").Value = 222
$ws.Range("G3").Value = 322
$ws.Range("H3").Value = 248
$ws.Range("I3").Value = 248
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3440
$ws.Range("L3").Value = 1352
$ws.Range("M3").Value = 2088
$ws.Range("N3").Value = 2083
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 210
$ws.Range("Q3").Value = 218
$ws.Range("R3").Value = -308
$ws.Range("S3").Value = 65
$ws.Range("T3").Value = 312
$ws.Range("U3").Value = -94
$ws.Range("V3").Value = 819
$ws.Range("W3").Value = 5.9
$ws.Range("X3").Value = 6.58
$ws.Range("Y3").Value = 12.61
$ws.Range("Z3").Value = 7.63
$ws.Range("AA3").Value = 64.73999999999999
$ws.Range("AB3").Value = 892.13
$ws.Range("AC3").Value = 590
$ws.Range("AD3").Value = 13.11
$ws.Range("AE3").Value = 4969
$ws.Range("AF3").Value = 1.56
$ws.Range("AG3").Value = 75
$ws.Range("AH3").Value = 0.97
$ws.Range("AI3").Value = 12.79
$ws.Range("AJ3").Value = 36720000

# Row 4
$ws.Range("D4").Value = 3779
$ws.Range("E4").Value = 156
$ws.Range("F4").Value = 156
$ws.Range("G4").Value = 558
$ws.Range("H4").Value = 405
$ws.Range("I4").Value = 405
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 3993
$ws.Range("L4").Value = 1530
$ws.Range("M4").Value = 2463
$ws.Range("N4").Value = 2455
$ws.Range("O4").Value = 8
$ws.Range("P4").Value = 210
$ws.Range("Q4").Value = 392
$ws.Range("R4").Value = 46
$ws.Range("S4").Value = -17
$ws.Range("T4").Value = 447
$ws.Range("U4").Value = -55
$ws.Range("V4").Value = 827
$ws.Range("W4").Value = 4.14
$ws.Range("X4").Value = 10.71
$ws.Range("Y4").Value = 17.87
$ws.Range("Z4").Value = 10.89
$ws.Range("AA4").Value = 62.11
$ws.Range("AB4").Value = 1070.16
$ws.Range("AC4").Value = 965
$ws.Range("AD4").Value = 6.67
$ws.Range("AE4").Value = 5854
$ws.Range("AF4").Value = 1.1
$ws.Range("AG4").Value = 80
$ws.Range("AH4").Value = 1.24
$ws.Range("AI4").Value = 8.34
$ws.Range("AJ4").Value = 36720000

# Row 5
$ws.Range("D5").Value = 3863
$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 167
$ws.Range("G5").Value = 186
$ws.Range("H5").Value = 118
$ws.Range("I5").Value = 118
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3804
$ws.Range("L5").Value = 1298
$ws.Range("M5").Value = 2506
$ws.Range("N5").Value = 2506
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 210
$ws.Range("Q5").Value = 294
$ws.Range("R5").Value = -144
$ws.Range("S5").Value = -95
$ws.Range("T5").Value = 495
$ws.Range("U5").Value = -201
$ws.Range("V5").Value = 778
$ws.Range("W5").Value = 4.32
$ws.Range("X5").Value = 3.05
$ws.Range("Y5").Value = 4.75
$ws.Range("Z5").Value = 3.02
$ws.Range("AA5").Value = 51.81
$ws.Range("AB5").Value = 1106.66
$ws.Range("AC5").Value = 280
$ws.Range("AD5").Value = 18.63
$ws.Range("AE5").Value = 5975
$ws.Range("AF5").Value = 0.87
$ws.Range("AG5").Value = 90
$ws.Range("AH5").Value = 1.72
$ws.Range("AI5").Value = 32.29
$ws.Range("AJ5").Value = 36720000

# Row 6
$ws.Range("D6").Value = 4149
$ws.Range("E6").Value = 149
$ws.Range("F6").Value = 149
$ws.Range("G6").Value = 156
$ws.Range("H6").Value = 118
$ws.Range("I6").Value = 118
$ws.Range("K6").Value = 3928
$ws.Range("L6").Value = 1318
$ws.Range("M6").Value = 2610
$ws.Range("N6").Value = 2610
$ws.Range("P6").Value = 210
$ws.Range("Q6").Value = 388
$ws.Range("R6").Value = -499
$ws.Range("S6").Value = -21
$ws.Range("T6").Value = 381
$ws.Range("U6").Value = 7
$ws.Range("V6").Value = 800
$ws.Range("W6").Value = 3.6
$ws.Range("X6").Value = 2.86
$ws.Range("Y6").Value = 4.63
$ws.Range("Z6").Value = 3.06
$ws.Range("AA6").Value = 50.51
$ws.Range("AB6").Value = 1150.86
$ws.Range("AC6").Value = 282
$ws.Range("AD6").Value = 13.35
$ws.Range("AE6").Value = 6225
$ws.Range("AF6").Value = 0.6
$ws.Range("AG6").Value = 95
$ws.Range("AH6").Value = 2.52
$ws.Range("AI6").Value = 33.84
$ws.Range("AJ6").Value = 36720000

# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
